$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update rows 3-11 (productId / score / timestamp revisions) ---
$updates = @(
    @{Row=3;  B="6728e8a8071b8fcf4f501df0"; C=0.7;   D="2025-02-12T11:31:56.846Z"},
    @{Row=4;  B="6728e9ab071b8fcf4f501df6"; C=0.88;  D="2025-04-13T11:31:56.846Z"},
    @{Row=5;  B="6728e9e6071b8fcf4f501dfc"; C=0.94;  D="2025-04-13T11:31:56.846Z"},
    @{Row=6;  B="6728ea62071b8fcf4f501e02"; C=0.82;  D="2025-04-13T11:31:56.846Z"},
    @{Row=7;  B="6743a3a8fd3ceed5b16a5e18"; C=0.82;  D="2025-04-13T11:31:56.846Z"},
    @{Row=8;  B="676137906c06138b1419f8a5"; C=0.94;  D="2025-04-13T11:31:56.846Z"},
    @{Row=9;  B="67f08cab1841d535b6af6f50"; C=0.76;  D="2025-04-13T11:31:56.846Z"},
    @{Row=10; B="67f08e231841d535b6af6f67"; C=0.82;  D="2025-04-13T11:31:56.846Z"},
    @{Row=11; B="6728eb4a071b8fcf4f501e0b"; C=0.76;  D="2025-05-23T04:46:31.245Z"}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 2).Value = $u.B
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 4).Value = $u.D
}

# --- Append new rows 48-54 for the new user 682326702fff19d415752f01 ---
$newRows = @(
    @{Row=48; A="682326702fff19d415752f01"; B="6728e9e6071b8fcf4f501dfc"; C=0.925; D="2025-05-23T04:41:14.096Z"},
    @{Row=49; A="682326702fff19d415752f01"; B="6728e93a071b8fcf4f501df3"; C=0.775; D="2025-05-23T09:19:25.598Z"},
    @{Row=50; A="682326702fff19d415752f01"; B="6728ea18071b8fcf4f501dff"; C=0.775; D="2025-05-23T09:19:36.004Z"},
    @{Row=51; A="682326702fff19d415752f01"; B="67f091181841d535b6af6f7b"; C=0.775; D="2025-05-23T09:19:59.350Z"},
    @{Row=52; A="682326702fff19d415752f01"; B="6728e9cd071b8fcf4f501df9"; C=0.925; D="2025-05-23T09:21:52.485Z"},
    @{Row=53; A="682326702fff19d415752f01"; B="6728f96acb86d3695fa1f4a6"; C=0.775; D="2025-05-23T09:21:00.865Z"},
    @{Row=54; A="682326702fff19d415752f01"; B="68067dd1286f80e4174d8736"; C=1;     D="2025-05-23T09:23:22.245Z"}
)

foreach ($n in $newRows) {
    $ws.Cells.Item($n.Row, 1).Value = $n.A
    $ws.Cells.Item($n.Row, 2).Value = $n.B
    $ws.Cells.Item($n.Row, 3).Value = $n.C
    $ws.Cells.Item($n.Row, 4).Value = $n.D
}
